$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "41.819.57"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "2.264.56"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "2.616.01"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "2.271.77"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").Value = "41.747.44"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  -5.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.106"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").Value = "2.006.59"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.93%  "
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "
